# Added test data for Switzerland market
#
# Duplicate the "Czech" sheet (it already has the right layout/formatting
# used by every market sheet in this workbook) to create the new "Swiss"
# sheet, then patch the market-specific cell values and view state.

$wb = $excel.ActiveWorkbook

$czech = $wb.Worksheets.Item("Czech")
$czech.Copy($null, $czech)

$swiss = $wb.Worksheets.Item(4)
$swiss.Name = "Swiss"
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2653"
$swiss.Rows.Item(11).RowHeight = 13.8
$swiss.Rows.Item(12).EntireRow.AutoFit()
$swiss.Range("A9").Select() | Out-Null

# Belgium tab keeps its old selection state, just moved to A8:A11
$belgium = $wb.Worksheets.Item("Belgium")
$belgium.Activate()
$belgium.Range("A8:A11").Select() | Out-Null

# Czech becomes the active/selected tab
$czech.Activate()
